$wb = $excel.ActiveWorkbook

# --- Update the conversion note text in Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.29 = 12688.17 pesos`n✅ 12688.17 pesos = 3.28 = 965.85 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures in the tasas sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 303.5
$wsTasas.Range("O10").Value = 3850.86
$wsTasas.Range("N12").Value = 3868.8
$wsTasas.Range("O12").Value = 294.5
